$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells from human-readable Indonesian labels to camelCase
# field names (matching the backend's normalized naming), and trim the
# trailing qualifier off the last header.
$ws.Range("A1").Value = "jenisKelamin"
$ws.Range("B1").Value = "organisasi"
$ws.Range("C1").Value = "ekstrakurikuler"
$ws.Range("D1").Value = "sertifikasiProfesi"
$ws.Range("E1").Value = "nilaiAkhir"
$ws.Range("F1").Value = "tempatMagang"
$ws.Range("G1").Value = "tempatKerja"
$ws.Range("H1").Value = "Durasi Mendapat Kerja"
